# Fix FileNotFoundError in save_input_files: the CONDUCTOR_files sheet
# listed dummy placeholder file names for the various EXTERNAL_* inputs
# (alphab_dummy.xlsx, bfield.xlsx, I_file_dummy.xlsx, flow_dummy.xlsx,
# Q_file_dummy.xlsx, strain_dummy.xlsx, spatial_discretization.xlsx).
# None of those files actually exist on disk, so replace each of them
# with the sentinel value "none" (already used elsewhere in the sheet
# to mean "no external file").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CONDUCTOR_files")

$ws.Range("E9").Value  = "none"   # EXTERNAL_ALPHAB          was alphab_dummy.xlsx
$ws.Range("E10").Value = "none"   # EXTERNAL_BFIELD          was bfield.xlsx
$ws.Range("E11").Value = "none"   # EXTERNAL_CURRENT         was I_file_dummy.xlsx
$ws.Range("E12").Value = "none"   # EXTERNAL_FLOW            was flow_dummy.xlsx
$ws.Range("E13").Value = "none"   # EXTERNAL_HEAT            was Q_file_dummy.xlsx
$ws.Range("E14").Value = "none"   # EXTERNAL_STRAIN          was strain_dummy.xlsx
$ws.Range("E15").Value = "none"   # EXTERNAL_GRID            was spatial_discretization.xlsx

# The workbook was re-saved with the CONDUCTOR_files sheet active/selected
# (instead of CONDUCTOR_input) and with cell E19 selected on it.
$ws.Activate()
$ws.Range("E19").Select()
